# HotFix adding missing QC fields to script output
# Inserts qc_reviewer_lanid / qc_notes / qc_status / qc_flags rows into the
# template_map Sheet1 table (and a missing conc_time_values "id" row),
# expanding the table from A1:C85 to A1:C102.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Insert the new blank rows, working from the bottom of the sheet
#        upward so that the row numbers used below for the still-untouched
#        upper rows stay valid while we work. ---
$ws.Rows("86:88").Insert()   # end of conc_time_values block (3 new rows)
$ws.Rows("79:82").Insert()   # end of series block (3 rows) + new conc_time_values "id" row
$ws.Rows("53:55").Insert()   # end of subjects block (3 rows)
$ws.Rows("41:43").Insert()   # end of studies block (3 rows)
$ws.Rows("14:17").Insert()   # end of documents block (4 rows)

# --- 2. Populate the new rows. Columns B/C are written first in the exact
#        order the new label strings first appear (qc_reviewer_lanid,
#        qc_notes, qc_status, qc_flags) so they match the expected
#        shared-string table ordering; column A is the section name, which
#        already exists in the shared strings. ---

# documents (rows 14-17)
$ws.Range("A14").Value = "documents"
$ws.Range("B14").Value = "qc_reviewer_lanid"
$ws.Range("C14").Value = "qc_reviewer_lanid"

$ws.Range("A15").Value = "documents"
$ws.Range("B15").Value = "qc_notes"
$ws.Range("C15").Value = "qc_notes"

$ws.Range("A16").Value = "documents"
$ws.Range("B16").Value = "qc_status"
$ws.Range("C16").Value = "qc_status"

$ws.Range("A17").Value = "documents"
$ws.Range("B17").Value = "qc_flags"
$ws.Range("C17").Value = "qc_flags"

# studies (rows 45-47)
$ws.Range("A45").Value = "studies"
$ws.Range("B45").Value = "qc_notes"
$ws.Range("C45").Value = "qc_notes"

$ws.Range("A46").Value = "studies"
$ws.Range("B46").Value = "qc_status"
$ws.Range("C46").Value = "qc_status"

$ws.Range("A47").Value = "studies"
$ws.Range("B47").Value = "qc_flags"
$ws.Range("C47").Value = "qc_flags"

# subjects (rows 60-62)
$ws.Range("A60").Value = "subjects"
$ws.Range("B60").Value = "qc_notes"
$ws.Range("C60").Value = "qc_notes"

$ws.Range("A61").Value = "subjects"
$ws.Range("B61").Value = "qc_status"
$ws.Range("C61").Value = "qc_status"

$ws.Range("A62").Value = "subjects"
$ws.Range("B62").Value = "qc_flags"
$ws.Range("C62").Value = "qc_flags"

# series (rows 89-91)
$ws.Range("A89").Value = "series"
$ws.Range("B89").Value = "qc_notes"
$ws.Range("C89").Value = "qc_notes"

$ws.Range("A90").Value = "series"
$ws.Range("B90").Value = "qc_status"
$ws.Range("C90").Value = "qc_status"

$ws.Range("A91").Value = "series"
$ws.Range("B91").Value = "qc_flags"
$ws.Range("C91").Value = "qc_flags"

# conc_time_values gains a new leading "id" row (row 92)
$ws.Range("A92").Value = "conc_time_values"
$ws.Range("B92").Value = "id"
$ws.Range("C92").Value = "id"

# conc_time_values (rows 100-102)
$ws.Range("A100").Value = "conc_time_values"
$ws.Range("B100").Value = "qc_notes"
$ws.Range("C100").Value = "qc_notes"

$ws.Range("A101").Value = "conc_time_values"
$ws.Range("B101").Value = "qc_status"
$ws.Range("C101").Value = "qc_status"

$ws.Range("A102").Value = "conc_time_values"
$ws.Range("B102").Value = "qc_flags"
$ws.Range("C102").Value = "qc_flags"

# --- 3. Refresh the AutoFilter range and the _FilterDatabase defined name
#        to match the new extent used by the authors (A1:C99), and update
#        the view so the new rows are visible/selected, matching how the
#        sheet was left after the edit. ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:C99").AutoFilter()

$ws.Range("B100:C102").Select()
$ws.Application.ActiveWindow.ScrollRow = 85
